# leetcode_everyday worksheet update:
#  - mark problem in row 25 (G column) as reviewed on 2021-12-16 (serial 44546)
#    instead of the "not yet reviewed" label
#  - append a new row (43) for "No72. Edit Distance"
#  - move the sheet's visible selection down to roughly where the new row is

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G25: switch from the "未复习" (not reviewed) placeholder to an actual
#     last-reviewed date -------------------------------------------------
$ws.Range("G25").Value = 44546

# --- new row 43: No72. Edit Distance ------------------------------------
# Copy the cell formatting (styles only) from row 40, which has the same
# shape (A:G populated, no H cell) as the row we are about to add.
$ws.Range("A40:G40").Copy()
$ws.Range("A43:G43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A43").Value = "No72. Edit Distance"
$ws.Range("B43").Value = "Hard"
$ws.Range("C43").Value = "https://leetcode.com/problems/edit-distance/"
$ws.Range("D43").Value = 44546
$ws.Range("E43").Value = "动态规划，字符串编辑"
$ws.Range("F43").Value = "和1143类似的dp数组设置，注意初始值设置"
$ws.Range("G43").Value = "未复习"

# Hyperlinks.Add() re-styles the target cell with the built-in "Hyperlink"
# look, so restore C43's real (copied-from-C40) formatting afterwards.
$ws.Hyperlinks.Add($ws.Range("C43"), "https://leetcode.com/problems/edit-distance/")
$ws.Range("C40").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows("43").RowHeight = 28

# --- move the selection/scroll position down towards the new row -------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F47").Select() | Out-Null
